$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H107").Value = 1514.6316
$ws.Range("I107").Value = 1487.6666
$ws.Range("K107").Value = 1487.6666
$ws.Range("M107").Value = 432.3334

$ws.Range("H137").Value = 1925.1538
$ws.Range("I137").Value = 1738.909
$ws.Range("J137").Value = 2949.5
$ws.Range("K137").Value = 5216.727000000001
$ws.Range("L137").Value = 8848.5
$ws.Range("M137").Value = -2666.727000000001
$ws.Range("N137").Value = -13948.5

$ws.Range("H138").Value = 5953.387
$ws.Range("J138").Value = 6702.6
$ws.Range("L138").Value = 20107.8
$ws.Range("N138").Value = -30387.8

$ws.Range("H141").Value = 1997.5
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").Value = $null

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 10640.371
$ws.Range("I32").Value = 9770.091
$ws.Range("K32").Value = 9770.091
$ws.Range("M32").Value = -9483.091

$ws.Range("H61").Value = 1427
$ws.Range("I61").Value = 1538.7142
$ws.Range("J61").Value = 1166.3334
$ws.Range("K61").Value = 1538.7142
$ws.Range("L61").Value = 1166.3334
$ws.Range("M61").Value = -1326.7142
$ws.Range("N61").Value = -1590.3334

$ws.Range("H132").Value = 3197.5
$ws.Range("J132").Value = 3965.6667
$ws.Range("L132").Value = 11897.0001
$ws.Range("N132").Value = -16957.0001

$ws.Range("H136").Value = 1427
$ws.Range("I136").Value = 1538.7142
$ws.Range("J136").Value = 1166.3334
$ws.Range("K136").Value = 4616.142599999999
$ws.Range("L136").Value = 3499.0002
$ws.Range("M136").Value = -2066.142599999999
$ws.Range("N136").Value = -8599.0002

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1189.1428
$ws.Range("I134").Value = 1189.1428
$ws.Range("K134").Value = 3567.4284
$ws.Range("M134").Value = -1032.4284

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1974.8889
$ws.Range("I31").Value = 1537.3846
$ws.Range("J31").Value = 3112.4
$ws.Range("K31").Value = 1537.3846
$ws.Range("L31").Value = 3112.4
$ws.Range("M31").Value = -1242.3846
$ws.Range("N31").Value = -3702.4

$ws.Range("H34").Value = 1974.8889
$ws.Range("I34").Value = 1537.3846
$ws.Range("J34").Value = 3112.4
$ws.Range("K34").Value = 1537.3846
$ws.Range("L34").Value = 3112.4
$ws.Range("M34").Value = -1335.3846
$ws.Range("N34").Value = -3516.4

$ws.Range("H58").Value = 6027.5
$ws.Range("I58").Value = 2990
$ws.Range("J58").Value = 7040
$ws.Range("K58").Value = 2990
$ws.Range("L58").Value = 7040
$ws.Range("M58").Value = -2787
$ws.Range("N58").Value = -7446

$ws.Range("H86").Value = 13328.8
$ws.Range("J86").Value = 14412.125
$ws.Range("L86").Value = 14412.125
$ws.Range("N86").Value = -16658.125

$ws.Range("H89").Value = 13328.8
$ws.Range("J89").Value = 14412.125
$ws.Range("L89").Value = 72060.625
$ws.Range("N89").Value = -83292.625

$ws.Range("H99").Value = 4500
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 4500
$ws.Range("K99").Value = 0
$ws.Range("L99").Value = 4500
$ws.Range("M99").Value = $null
$ws.Range("N99").Value = -7496

$ws.Range("H122").Value = 8223.25
$ws.Range("I122").Value = 7450
$ws.Range("K122").Value = 22350
$ws.Range("M122").Value = -19900

$ws.Range("H126").Value = 4500
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 4500
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 13500
$ws.Range("M126").Value = $null
$ws.Range("N126").Value = -18440

$ws.Range("H136").Value = 6027.5
$ws.Range("I136").Value = 2990
$ws.Range("J136").Value = 7040
$ws.Range("K136").Value = 8970
$ws.Range("L136").Value = 21120
$ws.Range("M136").Value = -6420
$ws.Range("N136").Value = -26220

$ws.Range("H141").Value = 698665
$ws.Range("J141").Value = 698665
$ws.Range("L141").Value = 698665
$ws.Range("N141").Value = -709025

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 95
$ws.Range("I4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("M4").Value = $null

$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 0
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 0
$ws.Range("M17").Value = $null
$ws.Range("N17").Value = $null

$ws.Range("H34").Value = 3499.75
$ws.Range("J34").Value = 6500
$ws.Range("L34").Value = 19500
$ws.Range("N34").Value = -19668

$ws.Range("H39").Value = 19996
$ws.Range("J39").Value = 19996
$ws.Range("L39").Value = 59988
$ws.Range("N39").Value = -60576

$ws.Range("H55").Value = 20000
$ws.Range("J55").Value = 20000
$ws.Range("L55").Value = 60000
$ws.Range("N55").Value = -60354

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H31").Value = 2366.6667
$ws.Range("I31").Value = 2366.6667
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 2366.6667
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = -2074.6667
$ws.Range("N31").Value = $null

$ws.Range("H37").Value = 2366.6667
$ws.Range("I37").Value = 2366.6667
$ws.Range("J37").Value = 0
$ws.Range("K37").Value = 2366.6667
$ws.Range("L37").Value = 0
$ws.Range("M37").Value = -2089.6667
$ws.Range("N37").Value = $null

$ws.Range("H122").Value = 3864.9
$ws.Range("I122").Value = 2781.6667
$ws.Range("J122").Value = 7114.6
$ws.Range("K122").Value = 8345.000100000001
$ws.Range("L122").Value = 21343.8
$ws.Range("M122").Value = -5895.000100000001
$ws.Range("N122").Value = -26243.8

$ws.Range("H132").Value = 5999
$ws.Range("I132").Value = 4001.6667
$ws.Range("J132").Value = 7497
$ws.Range("K132").Value = 12005.0001
$ws.Range("L132").Value = 22491
$ws.Range("M132").Value = -9475.000100000001
$ws.Range("N132").Value = -27551

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4097.2856
$ws.Range("I7").Value = 4097.2856
$ws.Range("K7").Value = 4097.2856
$ws.Range("M7").Value = -3985.2856

$ws.Range("H32").Value = 0
$ws.Range("I32").Value = 0
$ws.Range("K32").Value = 0
$ws.Range("M32").Value = $null

$ws.Range("H40").Value = 4749.5
$ws.Range("I40").Value = 4749.5
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 4749.5
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -4613.5
$ws.Range("N40").Value = $null

$ws.Range("H125").Value = 60000
$ws.Range("J125").Value = 60000
$ws.Range("L125").Value = 60000
$ws.Range("N125").Value = -69840

$ws.Range("H126").Value = 4097.2856
$ws.Range("I126").Value = 4097.2856
$ws.Range("K126").Value = 12291.8568
$ws.Range("M126").Value = -9821.856800000001

$ws.Range("H136").Value = 4438.1113
$ws.Range("I136").Value = 1927.7142
$ws.Range("J136").Value = 13224.5
$ws.Range("K136").Value = 5783.142599999999
$ws.Range("L136").Value = 39673.5
$ws.Range("M136").Value = -3233.142599999999
$ws.Range("N136").Value = -44773.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 4697.6875
$ws.Range("J81").Value = 6561.75
$ws.Range("L81").Value = 13123.5
$ws.Range("N81").Value = -15245.5

$ws.Range("H84").Value = 4697.6875
$ws.Range("J84").Value = 6561.75
$ws.Range("L84").Value = 65617.5
$ws.Range("N84").Value = -76225.5

$ws.Range("H126").Value = 1694.1177
$ws.Range("J126").Value = 1700
$ws.Range("L126").Value = 5100
$ws.Range("N126").Value = -10040

$ws.Range("H132").Value = 5307.8887
$ws.Range("I132").Value = 2018.25
$ws.Range("K132").Value = 6054.75
$ws.Range("M132").Value = -3524.75

$ws.Range("H136").Value = 1002.63635
$ws.Range("I136").Value = 1014.1
$ws.Range("J136").Value = 888
$ws.Range("K136").Value = 3042.3
$ws.Range("L136").Value = 2664
$ws.Range("M136").Value = -492.3000000000002
$ws.Range("N136").Value = -7764
